# Add kernel SVR parameters (kernel scale, epsilon, box constraint) to the
# pred_par.xlsx parameter sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header labels for the kernel SVR parameters (row 1)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Updated existing parameter values in row 2
$ws.Range("F2").Value = 0.001
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 12

# New kernel SVR parameter values (row 2)
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Match the view/selection state left behind after the edit
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("L5").Select()
